$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.358.73"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "3.898.13"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'529.57"
$ws.Range("E5").Value = "  +9.88%  "
$ws.Range("D6").Value = "'144.98"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  -1.47%  "
$ws.Range("D9").Value = "'0.719"
$ws.Range("E9").Value = "  -3.04%  "
$ws.Range("E10").Value = "  -2.37%  "
$ws.Range("E11").Value = "  -4.31%  "
$ws.Range("D12").Value = "'42.17"
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("D13").Value = "4.523.07"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").Value = "'10.27"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "3.924.95"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "'14.00"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("E18").Value = "  +7.15%  "
$ws.Range("D19").Value = "'19.83"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "69.288.77"
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").Value = "'425.14"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").Value = "'3.39"
$ws.Range("E22").Value = "  -5.25%  "
$ws.Range("E23").Value = "  -3.79%  "
$ws.Range("D24").Value = "'88.16"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").Value = "'4.05"
$ws.Range("E25").Value = "  +10.57%  "
$ws.Range("E26").Value = "  -7.34%  "
$ws.Range("D27").Value = "'10.59"
$ws.Range("E27").Value = "  -3.45%  "
$ws.Range("D28").Value = "'36.35"
$ws.Range("E28").Value = "  -1.93%  "
$ws.Range("D29").Value = "'690.43"
$ws.Range("E29").Value = "  -3.82%  "
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("E31").Value = "  -2.58%  "
$ws.Range("D32").Value = "'2.83"
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("D33").Value = "'68.99"
$ws.Range("E33").Value = "  +11.84%  "
$ws.Range("D34").Value = "0.0₃0878"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "'0.436"
$ws.Range("E35").Value = "  +9.28%  "
$ws.Range("D36").Value = "'5.93"
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("D37").Value = "'40.04"
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("E38").Value = "  +2.49%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").Value = "'3.26"
$ws.Range("E41").Value = "  +6.37%  "
$ws.Range("D42").Value = "'3.21"
$ws.Range("E42").Value = "  +9.02%  "
$ws.Range("D43").Value = "'0.0481"
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("D44").Value = "'2.81"
$ws.Range("E44").Value = "  -5.38%  "
$ws.Range("D45").Value = "'3.42"
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("D46").Value = "'0.000288"
$ws.Range("E46").Value = "  +15.74%  "
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("E48").Value = "  +6.95%  "
$ws.Range("D49").Value = "'145.93"
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("D50").Value = "2.746.33"
$ws.Range("E50").Value = "  +14.88%  "
$ws.Range("D51").Value = "0.0₆0343"
$ws.Range("E51").Value = "  -3.69%  "
